# Publish terminology IG 2.0.2 (#54)
# Updates the Metadata sheet: Version, Status, Experimental (cleared), Date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "1.8.2"
$ws.Range("B6").Value = "active"
$ws.Range("B7").ClearContents()

# B8 holds a date-looking string ("2025-11-18") that must stay a literal text
# value (not get auto-parsed into a date serial number). Enter it as a
# formula that evaluates to the text, then convert the cell to a plain
# value in place so it round-trips as a shared string instead of a date.
$ws.Range("B8").Formula = '="2025-11-18"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
